$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pin data table (B:C:D) was reworked for the new ADC slave-read wiring:
# existing Dx / Luz N labels gained "/<destino>" suffixes, D0.2 and D0.1 swapped
# their row order, two extra readout rows (luz 17 / luz18) were appended, and a
# brand-new a2..a14 address block (rows 37-49) was added below the table.
$rowNums = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49)
$colB    = @("pines Arduino", 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 3, 4, "a14", "a13", "a12", "a11", "a10", "a9", "a8", "a7", "a6", "a5", "a4", "a3", "a2")
$colC    = @("cable paralelo", 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 27, 28, 29, 30, 31, 32, 33, 34, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, $null)
$colD    = @("uso", "Boton 1", "Boton 5", "Boton 8", "Boton 7", "Boton 4", "Boton 6", "Boton 3", "Boton 2", "Dir 0", "Dir 1", "Dir 2", "Dir 3", "dis 1", "dis 2", "hd", "D0.0/Unidad1", "D0.3/Unidad4", "D0.2/Unidad3", "D0.1/unidad2", "D2/Decena3", "D1/Decena2", "D3/Decena4", "D0/Decena1", "Luz 12/Cuerpo 5", "Luz 11/Cuerpo 6", "Luz 3/Cuerpo 3", "Luz 5/Cuerpo 4", "Luz 7/M2", "Luz 8/Error", "Luz 16Cuerpo 2", "Luz 6/M1", "Luz 2/Cuerpo 1", "luz 17", "luz18", "h", "f", "e", "b", "g", "c", "d", "a", "k", "i", "l", "j", "adc")

for ($i = 0; $i -lt $rowNums.Count; $i++) {
    $r = $rowNums[$i]
    $ws.Range("B$r").Value = $colB[$i]
    if ($null -ne $colC[$i]) {
        $ws.Range("C$r").Value = $colC[$i]
    }
    $ws.Range("D$r").Value = $colD[$i]
}

# Column D got a bit wider to fit the longer "<pin>/<destino>" labels.
$ws.Columns.Item(4).ColumnWidth = 13.3

# Match the author's final selection/scroll position.
$ws.Range("E34").Select()
